# Actualización automática 2025-10-22 17:30:09
$wb = $excel.ActiveWorkbook

# --- Sheet: VENTAS POR GRUPO ---
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")
$ws1.Range("D41").Value = 783.55
$ws1.Range("L41").Value = 5892.92
$ws1.Range("M41").Value = 2264.12
$ws1.Range("D60").Value = "5 de 58"

# --- Sheet: VENTA MENSUAL ---
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")
$ws2.Range("F41").Value = 9199.790000000001
$ws2.Range("F60").Value = 37037.24000000001

# --- Sheet: CUMPLIMIENTO MENSUAL ---
$ws3 = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

$ws3.Range("D3").Value = 5110.66
$ws3.Range("E3").Value = 15276.8174217135
$ws3.Range("F3").Value = 0.2506764272149206

$ws3.Range("D11").Value = 10190.07
$ws3.Range("E11").Value = 9382.990249249699
$ws3.Range("F11").Value = 0.5206171068926546

$ws3.Range("D12").Value = 13020.19
$ws3.Range("E12").Value = 35603.87
$ws3.Range("F12").Value = 0.2677725800766123

$ws3.Range("D14").Value = 40544.3
$ws3.Range("E14").Value = 59353.69284188786
$ws3.Range("F14").Value = 0.4058570031949584
